$d = $word.ActiveDocument

# 1) Remove the stale "_GoBack" bookmark left after "color subtraction" from the
#    previous editing session.
$old = $d.Bookmarks("_GoBack")
$old.Delete()

# 2) Simulate the real edit: the author retyped the final "V" of "OpenCV" in the
#    VEXU bullet ("...object recognition using OpenCV"), which both leaves Word's
#    "last edit" _GoBack bookmark there and splits the run in two.
$rng = $d.Content
$found = $rng.Find.Execute("object recognition using OpenCV", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$endPos = $rng.End

# Drop the bookmark right before the final "V" - this marks the new _GoBack
# location and (as a side effect, matching how Word splits runs around
# bookmarks) separates "...OpenC" and "V" into distinct runs.
$pt = $d.Range($endPos - 1, $endPos - 1)
$d.Bookmarks.Add("_GoBack", $pt)
